# Serpina1e-Lrp1.xlsx update: refresh with new TPM-derived values.
#
# The previous data contained two "Sending cluster" blocks (ECs and FAPs),
# each with 3 rows (one per Target cluster: ECs, FAPs, MuSCs) -> 6 data rows.
# The refreshed data only keeps the "FAPs" sending-cluster block (3 rows),
# with updated (new TPM) numeric values throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "FAPs" sending-cluster block rows (previously rows 5-7); the
# remaining rows 2-4 (previously the "ECs" sending-cluster block) are
# overwritten below with the new data, becoming the sole surviving block.
$ws.Rows("5:7").Delete()

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Serpina1e"
$ws.Range("C2").Value = "Lrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 10.07790633333333
$ws.Range("H2").Value = 30.233719
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 20.42907845619711
$ws.Range("R2").Value = 183.861706105774
$ws.Range("S2").Value = 0.006596284565418616
$ws.Range("T2").Value = 0.006596284565418615

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Serpina1e"
$ws.Range("C3").Value = "Lrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 10.07790633333333
$ws.Range("H3").Value = 30.233719
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 2584.421688866067
$ws.Range("R3").Value = 23259.7951997946
$ws.Range("S3").Value = 0.8344762556643375
$ws.Range("T3").Value = 0.8344762556643374

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Serpina1e"
$ws.Range("C4").Value = "Lrp1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 10.07790633333333
$ws.Range("H4").Value = 30.233719
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 492.2076226837826
$ws.Range("R4").Value = 4429.868604154043
$ws.Range("S4").Value = 0.158927459770244
$ws.Range("T4").Value = 0.158927459770244
